# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The workbook records a serial-date value of 46075 in column C for every
# data row (rows 2-441); this should be bumped to 46076 (one day later)
# for all of them, leaving everything else untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46075) {
        $cell.Value2 = 46076
    }
}
